$p = $ppt.ActivePresentation

# Slide 1: Title "Header with inline code"
# Runs before: "Header" | " " | "with" | " " | "inline code"(Consolas)
# Runs after:  "Header " | "with " | "inline code"(Consolas)
$tr1 = $p.Slides.Item(1).Shapes.Item(1).TextFrame.TextRange
$tr1.Characters(1, 7).Text = "Header "
$tr1.Characters(8, 5).Text = "with "

# Slide 2: Title "Syntax highlighting"
# Runs before: "Syntax" | " " | "highlighting"
# Runs after:  "Syntax " | "highlighting"
$tr2 = $p.Slides.Item(2).Shapes.Item(1).TextFrame.TextRange
$tr2.Characters(1, 7).Text = "Syntax "

# Slide 3: Title "Two column slide"
# Runs before: "Two" | " " | "column" | " " | "slide"
# Runs after:  "Two " | "column " | "slide"
$tr3 = $p.Slides.Item(3).Shapes.Item(1).TextFrame.TextRange
$tr3.Characters(1, 4).Text = "Two "
$tr3.Characters(5, 7).Text = "column "
